$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change header from 'barcode' to 'product_id'
$ws.Range("A1").Value = "product_id"

# Row 2: barcode -> SKU
$ws.Range("A2").Value = "SKU123456"

# Row 3: 4607034370244 stays unchanged

# Row 4: barcode -> OZON SKU
$ws.Range("A4").Value = "OZON789012"

# Row 5: 9999999999999 stays unchanged

# Row 6: barcode -> UNKNOWN_SKU placeholder
$ws.Range("A6").Value = "UNKNOWN_SKU"
